# Form the consolidated report: update the "Absent" (column H) values
# for the rows whose attendance count changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
